$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 124, shifting existing rows 124:200 down to 125:201
$ws.Rows.Item(124).Insert()

# Fill the new row 124 with values, copying the constant categorical fields
# from the row below (which holds what used to be row 124's data) and
# setting the new data-specific values per the source update.
$ws.Range("A124").Value = 10
$ws.Range("B124").Value = "Vega Modelo de Temuco"
$ws.Range("C124").Value = "La Araucanía"
$ws.Range("D124").Value = 44582
$ws.Range("E124").Value = 9
$ws.Range("F124").Value = 100112039
$ws.Range("G124").Value = "Ciboulette"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 50
$ws.Range("K124").Value = 5000
$ws.Range("L124").Value = 5000
$ws.Range("M124").Value = 5000
$ws.Range("N124").Value = "$/docena de atados"
$ws.Range("O124").Value = "Provincia de Cautín"
$ws.Range("P124").Value = 1667
$ws.Range("Q124").Value = 3
$ws.Range("R124").Value = "Hortaliza"
